$d = $word.ActiveDocument

# Locate the "Simpla skisser till hemsidan 30/8" paragraph by its text rather than a
# hard-coded index, so the script keeps working even if earlier paragraphs shift.
function Find-ParagraphContaining($doc, $needle) {
    foreach ($p in $doc.Paragraphs) {
        if ($p.Range.Text -like "*$needle*") {
            return $p
        }
    }
    return $null
}

# --- The paragraphs involved in the edit -------------------------------
# "Simpla skisser till hemsidan 30/8"           -> the _GoBack bookmark moves to its end
# "Databasen ska vara klar senast den 31/8"     -> a new paragraph is inserted right before it
# "Inloggning ska fungera 2/9 "                 -> loses the _GoBack bookmark it used to carry

# 1) Remove the existing (hidden) _GoBack bookmark from the "Inloggning..." paragraph.
$oldBm = $d.Bookmarks.Item("_GoBack")
$oldBm.Delete()

# 2) Re-create _GoBack collapsed at the very end of the "Simpla skisser..." paragraph,
#    i.e. right after "30/8" and before the paragraph mark.
#    A zero-length Range placed exactly on the paragraph-mark boundary can't be targeted
#    directly, so nudge past it with a throwaway character, anchor the bookmark there,
#    then remove the throwaway character again - the collapsed bookmark stays put at
#    that same offset.
$skisser = Find-ParagraphContaining $d "Simpla skisser till hemsidan"
$endPos = $skisser.Range.End - 1
$d.Range($endPos, $endPos).InsertAfter("X")
$d.Bookmarks.Add("_GoBack", $d.Range($endPos, $endPos))
$d.Range($endPos, $endPos + 1).Delete()

# 3) Insert the new "PowerPoint presentation 31/8 14:30" paragraph right after
#    "Simpla skisser till hemsidan 30/8" (and thus right before "Databasen...").
$skisser = Find-ParagraphContaining $d "Simpla skisser till hemsidan"
$skisser.Range.InsertParagraphAfter()
$newPara = Find-ParagraphContaining $d "Simpla skisser till hemsidan"
$newPara.Next().Range.Text = "PowerPoint presentation 31/8 14:30"
